$d = $word.ActiveDocument

$d.Content.Find.Execute("Data Wrangling and Visualization", $true, $false, $false, $false, $false, $true, 1, $false, "Artificial Intelligence Policy", 2)

$d.Content.Find.Execute("F2025", $true, $false, $false, $false, $false, $true, 1, $false, "S2026", 2)

$d.Content.Find.Execute("Hello and welcome to the assignments page for our course. All weekly assignments will be posted here.", $true, $false, $false, $false, $false, $true, 1, $false, "Hello and welcome to the work page for our course. All assigned work for the course will be posted here, including debate guidelines and prompts, assignments, exam review guides, and project guidelines.", 2)

$d.Content.Find.Execute("Note that, while you can find assignments here, submission of assignments itself goes to in the", $true, $false, $false, $false, $false, $true, 1, $false, "Note that, while you can find assignments here, submission of assignments goes through the", 2)
